$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44270
$ws.Range("M2").Value = 85
$ws.Range('Q2').Value = '$/caja 14 kilos granel'
$ws.Range('R2').Value = 'Provincia del Elquí'
$ws.Range("S2").Value = 857
$ws.Range("T2").Value = 14

# Row 3
$ws.Range("D3").Value = 44323
$ws.Range('Q3').Value = '$/caja 14 kilos granel'

# Row 4
$ws.Range("D4").Value = 44242
$ws.Range("M4").Value = 45
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 12000
$ws.Range("P4").Value = 12000
$ws.Range('Q4').Value = '$/caja 15 kilos granel'
$ws.Range("S4").Value = 800
$ws.Range("T4").Value = 15

# Row 5
$ws.Range("D5").Value = 44260
$ws.Range("M5").Value = 56
$ws.Range("N5").Value = 13000
$ws.Range("O5").Value = 13000
$ws.Range("P5").Value = 13000
$ws.Range('Q5').Value = '$/caja 14 kilos empedrada'
$ws.Range('R5').Value = 'Provincia del Elquí'
$ws.Range("S5").Value = 929

# Row 6
$ws.Range("D6").Value = 44271
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 12000
$ws.Range('Q6').Value = '$/caja 14 kilos granel'
$ws.Range('R6').Value = 'Provincia del Elquí'
$ws.Range("S6").Value = 857
$ws.Range("T6").Value = 14

# Row 7
$ws.Range("D7").Value = 44239
$ws.Range("M7").Value = 70
$ws.Range("N7").Value = 15000
$ws.Range("O7").Value = 15000
$ws.Range("P7").Value = 15000
$ws.Range('Q7').Value = '$/caja 15 kilos granel'
$ws.Range('R7').Value = 'Provincia de Limarí'
$ws.Range("S7").Value = 1000
$ws.Range("T7").Value = 15

# Row 8
$ws.Range("D8").Value = 44312
$ws.Range("M8").Value = 68
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 14000
$ws.Range("P8").Value = 14000
$ws.Range('Q8').Value = '$/caja 14 kilos granel'
$ws.Range("S8").Value = 1000
$ws.Range("T8").Value = 14

# Row 9
$ws.Range("D9").Value = 44245
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 15000
$ws.Range('Q9').Value = '$/caja 15 kilos granel'
$ws.Range("T9").Value = 15

# Row 10
$ws.Range("D10").Value = 44320
$ws.Range("M10").Value = 45
$ws.Range("N10").Value = 14000
$ws.Range("O10").Value = 14000
$ws.Range("P10").Value = 14000
$ws.Range('Q10').Value = '$/caja 14 kilos granel'
$ws.Range("T10").Value = 14

# Row 11
$ws.Range("D11").Value = 44278
$ws.Range("N11").Value = 13000
$ws.Range("O11").Value = 13000
$ws.Range("P11").Value = 13000
$ws.Range('Q11').Value = '$/caja 14 kilos empedrada'
$ws.Range('R11').Value = 'Provincia del Elquí'
$ws.Range("S11").Value = 929

# Row 12
$ws.Range("D12").Value = 44252
$ws.Range("M12").Value = 60
$ws.Range('Q12').Value = '$/caja 14 kilos empedrada'

# Row 13
$ws.Range("D13").Value = 44259
$ws.Range("M13").Value = 80
$ws.Range("N13").Value = 12000
$ws.Range("O13").Value = 12000
$ws.Range("P13").Value = 12000
$ws.Range('Q13').Value = '$/caja 15 kilos empedrada'
$ws.Range('R13').Value = 'Provincia de Limarí'
$ws.Range("S13").Value = 800
$ws.Range("T13").Value = 15

# Row 14
$ws.Range("D14").Value = 44322
$ws.Range("M14").Value = 50

# Row 15
$ws.Range("D15").Value = 44316
$ws.Range("M15").Value = 48

# Row 16
$ws.Range("D16").Value = 44315
$ws.Range("M16").Value = 65
$ws.Range("N16").Value = 14000
$ws.Range("O16").Value = 14000
$ws.Range("P16").Value = 14000
$ws.Range('Q16').Value = '$/caja 14 kilos granel'
$ws.Range('R16').Value = 'Provincia de Limarí'
$ws.Range("S16").Value = 1000

# Row 17
$ws.Range("D17").Value = 44313
$ws.Range("M17").Value = 36
$ws.Range("N17").Value = 14000
$ws.Range("O17").Value = 14000
$ws.Range("P17").Value = 14000
$ws.Range('Q17').Value = '$/caja 14 kilos granel'
$ws.Range("T17").Value = 14

# Row 18
$ws.Range("D18").Value = 44314
$ws.Range("M18").Value = 56

# Row 19
$ws.Range("D19").Value = 44238
$ws.Range("M19").Value = 60
$ws.Range("N19").Value = 15000
$ws.Range("O19").Value = 15000
$ws.Range("P19").Value = 15000
$ws.Range('Q19').Value = '$/caja 15 kilos granel'
$ws.Range('R19').Value = 'Provincia de Limarí'
$ws.Range("S19").Value = 1000
$ws.Range("T19").Value = 15
